$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.103.46"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "2.349.95"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'544.81"
$ws.Range("E5").Value = "  +6.20%  "
$ws.Range("D6").Value = "'134.86"
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "2.347.78"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  +3.17%  "
$ws.Range("E13").Value = "  +6.79%  "
$ws.Range("D14").Value = "2.764.40"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").Value = "'23.58"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "58.066.79"
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "2.350.59"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("D19").Value = "'10.63"
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("D20").Value = "'335.47"
$ws.Range("E20").Value = "  +2.27%  "
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "'61.74"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'8.46"
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("E28").Value = "  +7.82%  "
$ws.Range("E29").Value = "  +5.15%  "
$ws.Range("D30").Value = "'170.32"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").Value = "'6.14"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  +17.83%  "
$ws.Range("D34").Value = "'18.48"
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E37").Value = "  +6.57%  "
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("E39").Value = "  +4.39%  "
$ws.Range("D40").Value = "'39.35"
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("D41").Value = "'148.74"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("D43").Value = "'286.82"
$ws.Range("E43").Value = "  +4.09%  "
$ws.Range("D44").Value = "'3.61"
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("D45").Value = "'19.30"
$ws.Range("E45").Value = "  +5.82%  "
$ws.Range("D46").Value = "'0.0926"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("D48").Value = "'0.563"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").Value = "'0.385"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("D50").Value = "'0.0218"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("E51").Value = "  +3.25%  "
